# Auto-generated edit script: updates Gungnir_Profits cost/profit columns
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 156250850
$ws.Range("J98").Value = 4500
$ws.Range("L98").Value = 4500
$ws.Range("N98").Value = -7496
$ws.Range("H122").Value = 156250850
$ws.Range("J122").Value = 4500
$ws.Range("L122").Value = 13500
$ws.Range("N122").Value = -18400
$ws.Range("H125").Value = 1340.6666
$ws.Range("I125").Value = 1154.6666
$ws.Range("J125").Value = 1526.6666
$ws.Range("K125").Value = 10391.9994
$ws.Range("L125").Value = 13739.9994
$ws.Range("M125").Value = -7931.999400000001
$ws.Range("N125").Value = -18659.9994
$ws.Range("H137").Value = 1618.64
$ws.Range("I137").Value = 1343.875
$ws.Range("J137").Value = 2107.111
$ws.Range("K137").Value = 4031.625
$ws.Range("L137").Value = 6321.333
$ws.Range("M137").Value = -1481.625
$ws.Range("N137").Value = -11421.333
$ws.Range("H138").Value = 3792.5618
$ws.Range("I138").Value = 2657.4736
$ws.Range("J138").Value = 4100.657
$ws.Range("K138").Value = 7972.4208
$ws.Range("L138").Value = 12301.971
$ws.Range("M138").Value = -2832.4208
$ws.Range("N138").Value = -22581.971
$ws.Range("H139").Value = 49436.273
$ws.Range("J139").Value = 49436.273
$ws.Range("L139").Value = 49436.273
$ws.Range("N139").Value = -59716.273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 38613.035
$ws.Range("I45").Value = 46915.453
$ws.Range("K45").Value = 46915.453
$ws.Range("M45").Value = -46538.453
$ws.Range("H61").Value = 1962.1364
$ws.Range("I61").Value = 2095.7896
$ws.Range("J61").Value = 1860.56
$ws.Range("K61").Value = 2095.7896
$ws.Range("L61").Value = 1860.56
$ws.Range("M61").Value = -1883.7896
$ws.Range("N61").Value = -2284.56
$ws.Range("H110").Value = 1449.6666
$ws.Range("I110").Value = 800.0833
$ws.Range("J110").Value = 2748.8333
$ws.Range("K110").Value = 800.0833
$ws.Range("L110").Value = 2748.8333
$ws.Range("M110").Value = 1244.9167
$ws.Range("N110").Value = -6838.8333
$ws.Range("H122").Value = 4720.9443
$ws.Range("J122").Value = 1866.3334
$ws.Range("L122").Value = 5599.0002
$ws.Range("N122").Value = -10499.0002
$ws.Range("H131").Value = 48466.25
$ws.Range("J131").Value = 48466.25
$ws.Range("L131").Value = 48466.25
$ws.Range("N131").Value = -58546.25
$ws.Range("H132").Value = 22730912
$ws.Range("I132").Value = 32259244
$ws.Range("K132").Value = 96777732
$ws.Range("M132").Value = -96775202
$ws.Range("H136").Value = 1962.1364
$ws.Range("I136").Value = 2095.7896
$ws.Range("J136").Value = 1860.56
$ws.Range("K136").Value = 6287.3688
$ws.Range("L136").Value = 5581.68
$ws.Range("M136").Value = -3737.3688
$ws.Range("N136").Value = -10681.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 83334216
$ws.Range("I107").Value = 125000720
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 125000720
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = -124998800
$ws.Range("N107").Value = -5040
$ws.Range("H134").Value = 4455655
$ws.Range("I134").Value = 13663.5
$ws.Range("J134").Value = 7416982.5
$ws.Range("K134").Value = 40990.5
$ws.Range("L134").Value = 22250947.5
$ws.Range("M134").Value = -38455.5
$ws.Range("N134").Value = -22256017.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 100012910
$ws.Range("I99").Value = 250026190
$ws.Range("J99").Value = 4050
$ws.Range("K99").Value = 250026190
$ws.Range("L99").Value = 4050
$ws.Range("M99").Value = -250024692
$ws.Range("N99").Value = -7046
$ws.Range("H126").Value = 100012910
$ws.Range("I126").Value = 250026190
$ws.Range("J126").Value = 4050
$ws.Range("K126").Value = 750078570
$ws.Range("L126").Value = 12150
$ws.Range("M126").Value = -750076100
$ws.Range("N126").Value = -17090
$ws.Range("H132").Value = 17551240
$ws.Range("I132").Value = 1243.1
$ws.Range("J132").Value = 37051236
$ws.Range("K132").Value = 3729.3
$ws.Range("L132").Value = 111153708
$ws.Range("M132").Value = -1199.3
$ws.Range("N132").Value = -111158768

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2888
$ws.Range("I3").Value = 2646.1177
$ws.Range("K3").Value = 7938.353099999999
$ws.Range("M3").Value = -7826.353099999999
$ws.Range("H68").Value = 71436550
$ws.Range("I68").Value = 166667260
$ws.Range("K68").Value = 500001780
$ws.Range("M68").Value = -500000969
$ws.Range("H71").Value = 71436550
$ws.Range("I71").Value = 166667260
$ws.Range("K71").Value = 1500005340
$ws.Range("M71").Value = -1500001284
$ws.Range("H110").Value = 300
$ws.Range("I110").Value = 300
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 900
$ws.Range("N110").ClearContents()
$ws.Range("M110").Value = 3190

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 168.625
$ws.Range("I2").Value = 24.833334
$ws.Range("J2").Value = 600
$ws.Range("K2").Value = 24.833334
$ws.Range("L2").Value = 600
$ws.Range("M2").Value = 88.16666599999999
$ws.Range("N2").Value = -826
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H43").Value = 3050
$ws.Range("I43").Value = 1733.3334
$ws.Range("K43").Value = 1733.3334
$ws.Range("M43").Value = -1582.3334
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H80").Value = 5266579
$ws.Range("I80").Value = 3938.7693
$ws.Range("J80").Value = 16668966
$ws.Range("K80").Value = 3938.7693
$ws.Range("L80").Value = 16668966
$ws.Range("M80").Value = -2940.7693
$ws.Range("N80").Value = -16670962
$ws.Range("H83").Value = 5266579
$ws.Range("I83").Value = 3938.7693
$ws.Range("J83").Value = 16668966
$ws.Range("K83").Value = 19693.8465
$ws.Range("L83").Value = 83344830
$ws.Range("M83").Value = -14701.8465
$ws.Range("N83").Value = -83354814
$ws.Range("H102").Value = 864.4231
$ws.Range("I102").Value = 703.5263
$ws.Range("J102").Value = 1301.1428
$ws.Range("K102").Value = 703.5263
$ws.Range("L102").Value = 1301.1428
$ws.Range("M102").Value = 918.4737
$ws.Range("N102").Value = -4545.1428
$ws.Range("H107").Value = 3560.1516
$ws.Range("I107").Value = 468.93332
$ws.Range("J107").Value = 6136.1665
$ws.Range("K107").Value = 468.93332
$ws.Range("L107").Value = 6136.1665
$ws.Range("M107").Value = 1451.06668
$ws.Range("N107").Value = -9976.166499999999
$ws.Range("H126").Value = 4914.615
$ws.Range("I126").Value = 5100
$ws.Range("J126").Value = 4832.222
$ws.Range("K126").Value = 15300
$ws.Range("L126").Value = 14496.666
$ws.Range("M126").Value = -12830
$ws.Range("N126").Value = -19436.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2421.5557
$ws.Range("I7").Value = 2349.25
$ws.Range("K7").Value = 2349.25
$ws.Range("M7").Value = -2237.25
$ws.Range("H40").Value = 31252188
$ws.Range("I40").Value = 2000.8
$ws.Range("J40").Value = 83335830
$ws.Range("K40").Value = 2000.8
$ws.Range("L40").Value = 83335830
$ws.Range("M40").Value = -1864.8
$ws.Range("N40").Value = -83336102
$ws.Range("H46").Value = 4630224
$ws.Range("I46").Value = 6945010.5
$ws.Range("J46").Value = 650
$ws.Range("K46").Value = 6945010.5
$ws.Range("L46").Value = 650
$ws.Range("M46").Value = -6944822.5
$ws.Range("N46").Value = -1026
$ws.Range("H122").Value = 21147.268
$ws.Range("I122").Value = 25563.545
$ws.Range("K122").Value = 76690.63499999999
$ws.Range("M122").Value = -74240.63499999999
$ws.Range("H126").Value = 2421.5557
$ws.Range("I126").Value = 2349.25
$ws.Range("K126").Value = 7047.75
$ws.Range("M126").Value = -4577.75
$ws.Range("H132").Value = 21745698
$ws.Range("I132").Value = 100003310
$ws.Range("J132").Value = 7471.6943
$ws.Range("K132").Value = 300009930
$ws.Range("L132").Value = 22415.0829
$ws.Range("M132").Value = -300007400
$ws.Range("N132").Value = -27475.0829

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 18204308
$ws.Range("I132").Value = 32293116
$ws.Range("J132").Value = 6263.2085
$ws.Range("K132").Value = 96879348
$ws.Range("L132").Value = 18789.6255
$ws.Range("M132").Value = -96876818
$ws.Range("N132").Value = -23849.6255
$ws.Range("H136").Value = 3593.2292
$ws.Range("I136").Value = 7407.1763
$ws.Range("J136").Value = 1501.7097
$ws.Range("K136").Value = 22221.5289
$ws.Range("L136").Value = 4505.1291
$ws.Range("M136").Value = -19671.5289
$ws.Range("N136").Value = -9605.1291
